# Fall 2022 Week 10 ("day-after") inputs: fill in actual results for the
# column K matches that were previously marked "A" (Available/TBD).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 (rows 3-10) ---
$ws.Range("K3").Value  = "L"    # Daniel Burcham
$ws.Range("K4").Value  = "DNP"  # Leo Hayward
$ws.Range("K5").Value  = "L"    # Laura Thompson
$ws.Range("K6").Value  = "DNP"  # Kim Quan
$ws.Range("K7").Value  = "L"    # Scott Berry
$ws.Range("K8").Value  = "DNP"  # Jason Liess
$ws.Range("K9").Value  = "W"    # Jason Bohrer
$ws.Range("K10").Value = "W"    # Dan Aquino

# --- Table 2 (rows 15-22) ---
$ws.Range("K15").Value = "L"    # Jason Bohrer
$ws.Range("K16").Value = "L"    # Jason Liess
$ws.Range("K17").Value = "DNP"  # Daniel Burcham
$ws.Range("K18").Value = "W"    # Scott Berry
$ws.Range("K19").Value = "L"    # Dan Aquino
$ws.Range("K20").Value = "W"    # Ashley Daniels
$ws.Range("K21").Value = "DNP"  # Adrian Warden
$ws.Range("K22").Value = "DNP"  # Shelia Lowe

# Move the live selection to match the author's final cursor position.
$ws.Range("K25").Select()
